# feat: add 2022-Q4 data
$wb = $excel.ActiveWorkbook

# --- Step 1: update the "总计" (summary) sheet ---
# The existing row 2 (2021-Q4) moves down to row 3, and row 2 becomes the
# new 2022-Q4 entry.
$wsTotal = $wb.Worksheets.Item(1)

# Copy the old row-2 data/format down into row 3 first (this preserves the
# style used on column A).
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = $wsTotal.Range("B2").Value2
$wsTotal.Range("C3").Value = $wsTotal.Range("C2").Value2
$wsTotal.Range("D3").Value = $wsTotal.Range("D2").Value2

# Now overwrite row 2 with the new 2022-Q4 values.
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 4
$wsTotal.Range("D2").Value = 2.87

# --- Step 2: insert a new worksheet "2022-Q4" right after the "总计" sheet,
# so the final sheet order is 总计, 2022-Q4, 2021-Q4. ---
$wsQ4_2022 = $wb.Worksheets.Add($null, $wsTotal)
$wsQ4_2022.Name = "2022-Q4"

# Apply the header-row style (B1:H1) and the row-label style (A2:A5) used
# throughout the rest of the workbook before filling in values.
$wsTotal.Range("B1").Copy()
$wsQ4_2022.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ4_2022.Range("A2:A5").PasteSpecial(-4122)

# --- Step 3: populate the new "2022-Q4" sheet ---
$wsQ4_2022.Range("B1").Value = "基金代码"
$wsQ4_2022.Range("C1").Value = "基金名称"
$wsQ4_2022.Range("D1").Value = "基金规模"
$wsQ4_2022.Range("E1").Value = "股票总仓位"
$wsQ4_2022.Range("F1").Value = "仓位占比"
$wsQ4_2022.Range("G1").Value = "持有市值(亿元)"
$wsQ4_2022.Range("H1").Value = "仓位排名"

# Columns B:G (rows 2-5) hold numeric-looking codes/figures that must stay
# text, so force a text format before assigning them.
$txtRange = $wsQ4_2022.Range("B2:G5")
$txtRange.NumberFormat = "@"

$wsQ4_2022.Range("A2").Value = 0
$wsQ4_2022.Range("B2").Value = "968029"
$wsQ4_2022.Range("C2").Value = "恒生指数基金M类人民币（对冲）份额"
$wsQ4_2022.Range("D2").Value = "27.13"
$wsQ4_2022.Range("E2").Value = "99.24"
$wsQ4_2022.Range("F2").Value = "8.32"
$wsQ4_2022.Range("G2").Value = "2.2572"
$wsQ4_2022.Range("H2").Value = 1

$wsQ4_2022.Range("A3").Value = 1
$wsQ4_2022.Range("B3").Value = "486001"
$wsQ4_2022.Range("C3").Value = "工银瑞信中国机会全球配置股票（QDII）人民币"
$wsQ4_2022.Range("D3").Value = "6.42"
$wsQ4_2022.Range("E3").Value = "93.86"
$wsQ4_2022.Range("F3").Value = "3.19"
$wsQ4_2022.Range("G3").Value = "0.2048"
$wsQ4_2022.Range("H3").Value = 3

$wsQ4_2022.Range("A4").Value = 2
$wsQ4_2022.Range("B4").Value = "009562"
$wsQ4_2022.Range("C4").Value = "工银全球股票（QDII）美元"
$wsQ4_2022.Range("D4").Value = "6.42"
$wsQ4_2022.Range("E4").Value = "93.86"
$wsQ4_2022.Range("F4").Value = "3.19"
$wsQ4_2022.Range("G4").Value = "0.2048"
$wsQ4_2022.Range("H4").Value = 3

$wsQ4_2022.Range("A5").Value = 3
$wsQ4_2022.Range("B5").Value = "009563"
$wsQ4_2022.Range("C5").Value = "工银全球股票（QDII）港币"
$wsQ4_2022.Range("D5").Value = "6.42"
$wsQ4_2022.Range("E5").Value = "93.86"
$wsQ4_2022.Range("F5").Value = "3.19"
$wsQ4_2022.Range("G5").Value = "0.2048"
$wsQ4_2022.Range("H5").Value = 3

# Remove the temporary number-format styling so the text cells end up with
# no explicit style, matching the plain text cells used elsewhere.
$txtRange.ClearFormats()

Write-Host "Done"
